$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 87 currently carries the "last row" date-only format (YYYY-MM-DD).
# Since a new row is being appended, row 87 becomes a regular data row,
# so give it the standard datetime format used by all the other rows.
$ws.Range("A87").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data as row 88.
$ws.Range("A88").Value = 45675
$ws.Range("B88").Value = 208
$ws.Range("C88").Value = 204
$ws.Range("D88").Value = 205

# The newest row takes on the special "last row" date-only format.
$ws.Range("A88").NumberFormat = "YYYY-MM-DD"
